$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 11 de Mayo de 2020 a las 23:35"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1381342
$ws.Range("C4").Value = 13704
$ws.Range("D4").Value = 260188
$ws.Range("E4").Value = 1039609
$ws.Range("F4").Value = 16451
$ws.Range("G4").Value = 758
$ws.Range("H4").Value = 81545

# Alemania (row 10)
$ws.Range("B10").Value = 172517
$ws.Range("C10").Value = 638
$ws.Range("E10").Value = 19264
$ws.Range("G10").Value = 84
$ws.Range("H10").Value = 7653

# Brasil (row 11)
$ws.Range("B11").Value = 166162
$ws.Range("C11").Value = 3463
$ws.Range("E11").Value = 89862
$ws.Range("G11").Value = 220
$ws.Range("H11").Value = 11343

# Canada (row 16)
$ws.Range("B16").Value = 69911
$ws.Range("C16").Value = 1063
$ws.Range("D16").Value = 32664
$ws.Range("E16").Value = 32255

# Reorder Etiopia / Cabo Verde: Cabo Verde now appears before Etiopia in the
# country list, with Cabo Verde's row (140) holding the fresh totals and
# Etiopia's row (141) retaining what used to be its own totals.
$ws.Range("A140").Value = "Cabo Verde"
$ws.Range("B140").Value = 260
$ws.Range("C140").Value = 14
$ws.Range("D140").Value = 58
$ws.Range("E140").Value = 200
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 2

$ws.Range("A141").Value = "Etiopia"
$ws.Range("B141").Value = 250
$ws.Range("C141").Value = 11
$ws.Range("D141").Value = 105
$ws.Range("E141").Value = 140
$ws.Range("F141").Value = 1
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 5
